$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sensor raw data values (columns B, C, D for rows 2-9) ---
$ws.Range("B2").Value = 645
$ws.Range("C2").Value = 665
$ws.Range("D2").Value = 642

$ws.Range("B3").Value = 415
$ws.Range("C3").Value = 425
$ws.Range("D3").Value = 410

$ws.Range("B4").Value = 305
$ws.Range("C4").Value = 330
$ws.Range("D4").Value = 300

$ws.Range("B5").Value = 245
$ws.Range("C5").Value = 260
$ws.Range("D5").Value = 235

$ws.Range("B6").Value = 200
$ws.Range("C6").Value = 180
$ws.Range("D6").Value = 199

$ws.Range("B7").Value = 165
$ws.Range("C7").Value = 160
$ws.Range("D7").Value = 175

$ws.Range("B8").Value = 135
$ws.Range("C8").Value = 130
$ws.Range("D8").Value = 175

$ws.Range("B9").Value = 110
$ws.Range("C9").Value = 110
$ws.Range("D9").Value = 180

# --- Highlight the newly-confirmed / re-tested readings in red ---
$ws.Range("C7:D9").Interior.Color = 255

# --- Update selection to reflect where attention was during the test ---
[void]$ws.Range("C7:C9").Select()
